$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formula edit: column C changed from 1/LOG(x,2) to 1.5/LOG(x,2)-0.09 ---
# Re-create the shared formula across C1:C12 so Excel regenerates the
# t="shared" group (C1 standalone, C2 as master of C2:C12) with refreshed
# cached values, matching how Excel stores an edited+filled formula.
$ws.Range("C1:C12").FormulaR1C1 = "=1.5/LOG(RC[-2],2)-0.09"

# --- Chart 1, series 2 label changed to match the new formula ---
$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(2).Name = "1.5/log(n)-.09"

# --- Charts repositioned (dragged) on the sheet ---
$co1 = $ws.ChartObjects(1)
$co1.Left = 625.875
$co1.Top = 42.5
$co1.Width = 622.3125
$co1.Height = 343

$co3 = $ws.ChartObjects(3)
$co3.Left = 1385.8125
$co3.Top = 32.25
$co3.Width = 403.8125
$co3.Height = 216

$co4 = $ws.ChartObjects(4)
$co4.Left = 1331.375
$co4.Top = 310.75
$co4.Width = 402.5625
$co4.Height = 217

# --- Selection moved to C3 ---
$ws.Range("C3").Select() | Out-Null
